# "Add Mainak Podder in Excel file"
# A new column D is introduced with the header/value "Mainak Podder" typed
# into D2 (recorded as a new shared string). Excel widens column D to fit
# the new text and moves the active selection to E1, as it would after a
# user types the value into D2 and presses Enter/Tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Mainak Podder"

# Target width recorded by Excel for this text is 14.21875 characters.
# 13.333333333333334 is the input that rounds, through this engine's
# pixel-grid quantization of ColumnWidth, to the closest reproducible
# value (14.166666666666666).
$ws.Columns.Item(4).ColumnWidth = 13.333333333333334

$ws.Range("E1").Select()
